$wb = $excel.ActiveWorkbook

# --- Update the "Mixed 0.5" / "Dm05" results row (row 6) on both sheets ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C6").Value = 29.629629629629601
$ws1.Range("D6").Value = 61.1111111111111
$ws1.Range("E6").Value = 39.072039072038997
$ws1.Range("F6").Value = 0.047577540079752602
$ws1.Range("G6").Value = 4.2496791680653798

$ws2 = $wb.Worksheets.Item("Sheet1 (2)")
$ws2.Range("C6").Value = 29.629629629629601
$ws2.Range("D6").Value = 61.1111111111111
$ws2.Range("E6").Value = 39.072039072038997
$ws2.Range("F6").Value = 0.047577540079752602
$ws2.Range("G6").Value = 4.2496791680653798

# --- Update selections on each sheet, then make "Sheet1 (2)" the active tab ---
[void]$ws1.Activate()
[void]$ws1.Range("C6:G6").Select()

[void]$ws2.Activate()
[void]$ws2.Range("P16").Select()
